$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.942.74"
$ws.Range("E2").Value = "  +3.78%  "
$ws.Range("D3").Value = "3.813.56"
$ws.Range("E3").Value = "  +4.84%  "
$ws.Range("E4").Value = "  -0.53%  "
$ws.Range("D5").Value = "'422.20"
$ws.Range("E5").Value = "  +4.59%  "
$ws.Range("D6").Value = "'129.70"
$ws.Range("E6").Value = "  -1.16%  "
$ws.Range("D7").Value = "3.815.43"
$ws.Range("E7").Value = "  +5.01%  "
$ws.Range("D8").Value = "'0.606"
$ws.Range("E8").Value = "  -1.64%  "
$ws.Range("D9").Value = "'0.998"
$ws.Range("D10").Value = "'0.721"
$ws.Range("E10").Value = "  +0.45%  "
$ws.Range("D11").Value = "'0.161"
$ws.Range("E11").Value = "  +1.70%  "
$ws.Range("D12").Value = "'0.0000349"
$ws.Range("E12").Value = "  +19.74%  "
$ws.Range("D13").Value = "'40.54"
$ws.Range("E13").Value = "  -2.31%  "
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").Value = "'10.17"
$ws.Range("E14").Value = "  +3.75%  "
$ws.Range("B15").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C15").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D15").Value = "4.404.35"
$ws.Range("E15").Value = "  +4.10%  "
$ws.Range("E16").Value = "  +16.17%  "
$ws.Range("B17").Value = "TRON"
$ws.Range("C17").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D17").Value = "'0.137"
$ws.Range("E17").Value = "  -0.43%  "
$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").Value = "3.812.68"
$ws.Range("E18").Value = "  +5.03%  "
$ws.Range("D19").Value = "'19.65"
$ws.Range("E19").Value = "  -0.30%  "
$ws.Range("D20").Value = "67.024.75"
$ws.Range("E20").Value = "  +3.81%  "
$ws.Range("D21").Value = "'1.09"
$ws.Range("E21").Value = "  +1.56%  "
$ws.Range("D22").Value = "'406.66"
$ws.Range("E22").Value = "  -2.33%  "
$ws.Range("D23").Value = "'15.55"
$ws.Range("E23").Value = "  +2.13%  "
$ws.Range("D24").Value = "'83.87"
$ws.Range("E24").Value = "  -1.56%  "
$ws.Range("D25").Value = "'3.04"
$ws.Range("E25").Value = "  +2.54%  "
$ws.Range("D26").Value = "'37.13"
$ws.Range("E26").Value = "  +4.62%  "
$ws.Range("D27").Value = "'10.00"
$ws.Range("E27").Value = "  +7.72%  "
$ws.Range("B28").Value = "PancakeSwap"
$ws.Range("C28").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D28").Value = "'3.21"
$ws.Range("E28").Value = "  +2.38%  "
$ws.Range("B29").Value = "LEO"
$ws.Range("C29").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D29").Value = "'5.46"
$ws.Range("E29").Value = "  +7.08%  "
$ws.Range("D30").Value = "'9.00"
$ws.Range("E30").Value = "  +30.35%  "
$ws.Range("D31").Value = "'717.88"
$ws.Range("E31").Value = "  +7.78%  "
$ws.Range("D32").Value = "'12.61"
$ws.Range("E32").Value = "  +1.62%  "
$ws.Range("E33").Value = "  +3.03%  "
$ws.Range("D34").Value = "'0.121"
$ws.Range("E34").Value = "  +3.44%  "
$ws.Range("D35").Value = "'0.999"
$ws.Range("E35").Value = "  +0.00%  "
$ws.Range("D36").Value = "'0.152"
$ws.Range("D37").Value = "'38.62"
$ws.Range("E37").Value = "  -3.78%  "
$ws.Range("D38").Value = "'55.15"
$ws.Range("E38").Value = "  -1.02%  "
$ws.Range("D39").Value = "'5.40"
$ws.Range("E39").Value = "  +24.47%  "
$ws.Range("D40").Value = "0.0₃0755"
$ws.Range("E40").Value = "  +22.95%  "
$ws.Range("D41").Value = "'0.0453"
$ws.Range("E41").Value = "  -1.05%  "
$ws.Range("E42").Value = "  +2.00%  "
$ws.Range("E43").Value = "  +0.17%  "
$ws.Range("D44").Value = "'0.134"
$ws.Range("E44").Value = "  -3.00%  "
$ws.Range("D45").Value = "'3.32"
$ws.Range("E45").Value = "  +0.98%  "
$ws.Range("D46").Value = "'143.58"
$ws.Range("E46").Value = "  +1.31%  "
$ws.Range("E47").Value = "  +1.42%  "
$ws.Range("B48").Value = "ARBITRUM"
$ws.Range("C48").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D48").Value = "'2.04"
$ws.Range("E48").Value = "  -0.37%  "
$ws.Range("B49").Value = "TheGraph"
$ws.Range("C49").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D49").Value = "'0.308"
$ws.Range("E49").Value = "  +7.26%  "
$ws.Range("D50").Value = "'2.81"
$ws.Range("E50").Value = "  +1.20%  "
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").Value = "'25.65"
$ws.Range("E51").Value = "  -3.46%  "

$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").Style = "Normal"
$ws.Range("D17").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").Style = "Normal"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").Style = "Normal"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").Style = "Normal"
